$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to text so numeric-looking strings (e.g. "1.000",
    # "0.9997") are not silently reinterpreted as numbers by Excel, then
    # restore the "Normal" style so no stray cell format/style is left
    # behind (keeps the workbook's style table identical to the source).
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Price (D) / Volume 1h (E) updates ---------------------------------
Set-TextValue "D2" "28.426.29"
Set-TextValue "E2" "  +4.30%  "

Set-TextValue "D3" "1.797.40"
Set-TextValue "E3" "  +1.27%  "

Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.09%  "

Set-TextValue "D5" "314.32"
Set-TextValue "E5" "  +0.41%  "

Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.06%  "

Set-TextValue "D7" "0.5463"
Set-TextValue "E7" "  +4.97%  "

Set-TextValue "D8" "0.3829"
Set-TextValue "E8" "  +3.95%  "

Set-TextValue "D9" "0.07611"

Set-TextValue "D10" "42.53"
Set-TextValue "E10" "  -0.57%  "

Set-TextValue "D11" "1.124"
Set-TextValue "E11" "  +3.21%  "

Set-TextValue "D12" "1.000"
Set-TextValue "E12" "  -0.10%  "

Set-TextValue "D13" "21.18"
Set-TextValue "E13" "  +3.37%  "

Set-TextValue "D14" "6.206"
Set-TextValue "E14" "  +2.21%  "

Set-TextValue "D15" "7.409"
Set-TextValue "E15" "  +6.49%  "

Set-TextValue "D16" "1.795.81"
Set-TextValue "E16" "  +1.40%  "

Set-TextValue "D17" "91.63"
Set-TextValue "E17" "  +3.03%  "

Set-TextValue "D18" "0.00001074"
Set-TextValue "E18" "  +2.55%  "

Set-TextValue "D19" "0.06458"
Set-TextValue "E19" "  +0.19%  "

Set-TextValue "D20" "0.9997"
Set-TextValue "E20" "  -0.08%  "

Set-TextValue "D21" "17.36"
Set-TextValue "E21" "  +3.73%  "

Set-TextValue "D22" "5.973"
Set-TextValue "E22" "  +2.68%  "

Set-TextValue "D23" "28.427.46"
Set-TextValue "E23" "  +4.14%  "

Set-TextValue "D24" "11.45"
Set-TextValue "E24" "  +1.84%  "

Set-TextValue "D25" "2.121"
Set-TextValue "E25" "  -0.09%  "

Set-TextValue "D26" "159.15"
Set-TextValue "E26" "  +2.60%  "

Set-TextValue "D27" "20.73"
Set-TextValue "E27" "  +2.69%  "

Set-TextValue "D28" "2.403"
Set-TextValue "E28" "  +3.44%  "

Set-TextValue "D29" "2.002.75"
Set-TextValue "E29" "  +1.32%  "

Set-TextValue "D30" "123.36"
Set-TextValue "E30" "  +1.75%  "

Set-TextValue "E31" "  +6.16%  "

Set-TextValue "D32" "0.1028"
Set-TextValue "E32" "  +5.02%  "

Set-TextValue "D33" "5.769"
Set-TextValue "E33" "  +3.68%  "

Set-TextValue "D34" "3.675"
Set-TextValue "E34" "  +1.52%  "

# --- Rows 35/36 swapped places (Algorand <-> Hedera) --------------------
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D35" "0.06766"
Set-TextValue "E35" "  +13.32%  "

$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D36" "0.2332"
Set-TextValue "E36" "  +15.50%  "

Set-TextValue "D37" "0.02324"
Set-TextValue "E37" "  +3.73%  "

Set-TextValue "D38" "5.171"
Set-TextValue "E38" "  +6.94%  "

Set-TextValue "D39" "8.789"
Set-TextValue "E39" "  +8.81%  "

Set-TextValue "D40" "11.72"
Set-TextValue "E40" "  +4.27%  "

Set-TextValue "D41" "0.6404"
Set-TextValue "E41" "  +4.36%  "

Set-TextValue "D42" "0.9994"
Set-TextValue "E42" "  -0.06%  "

Set-TextValue "D43" "1.160"
Set-TextValue "E43" "  +1.73%  "

Set-TextValue "D44" "1.403"
Set-TextValue "E44" "  -2.07%  "

Set-TextValue "D45" "13.65"
Set-TextValue "E45" "  +3.82%  "

Set-TextValue "D46" "0.5978"
Set-TextValue "E46" "  +3.75%  "

Set-TextValue "D47" "3.678"
Set-TextValue "E47" "  +1.41%  "

Set-TextValue "D48" "126.58"
Set-TextValue "E48" "  +4.46%  "

Set-TextValue "D49" "2.001"
Set-TextValue "E49" "  +6.20%  "

Set-TextValue "D50" "1.152"
Set-TextValue "E50" "  +3.30%  "

Set-TextValue "D51" "0.06932"
Set-TextValue "E51" "  +3.35%  "
